$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldText = "Allows users to visualize and share data (sensors values, media, web links etc.)"

$run1 = "Allows users to visualize and "
$run2 = "anonymously share data and digital content "
$run3 = "(sensors values, media, web links etc.)"

$full = $tr.Text
$idx0 = $full.IndexOf($oldText)
$start = $idx0 + 1

# Replace the tail first (run2+run3) while the run1 prefix is left untouched,
# so the still-to-be-used $start offset for the prefix stays valid even though
# the replacement text is longer than the original.
$oldTail = $oldText.Substring($run1.Length)
$tailStart = $start + $run1.Length
$tailRange = $tr.Characters($tailStart, $oldTail.Length)
$tailRange.Text = $run2 + $run3

# Now split the newly inserted tail into its own two runs (run2 / run3) by
# recomputing the offset of run3 from the just-written text.
$run3Start = $tailStart + $run2.Length
$run3Range = $tr.Characters($run3Start, $run3.Length)
$run3Range.Text = $run3
